$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 20859.572
$ws.Range("I21").Value = 14754.25
$ws.Range("K21").Value = 14754.25
$ws.Range("M21").Value = -14286.25
$ws.Range("H23").Value = 20859.572
$ws.Range("I23").Value = 14754.25
$ws.Range("K23").Value = 14754.25
$ws.Range("M23").Value = -14520.25
$ws.Range("H40").Value = 1720
$ws.Range("I40").Value = 2500
$ws.Range("J40").Value = 1525
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 1525
$ws.Range("M40").Value = -2325
$ws.Range("N40").Value = -1875
$ws.Range("H62").Value = 4065.625
$ws.Range("I62").Value = 3216.1765
$ws.Range("J62").Value = 6128.5713
$ws.Range("K62").Value = 3216.1765
$ws.Range("L62").Value = 6128.5713
$ws.Range("M62").Value = -2592.1765
$ws.Range("N62").Value = -7376.5713
$ws.Range("H65").Value = 4065.625
$ws.Range("I65").Value = 3216.1765
$ws.Range("J65").Value = 6128.5713
$ws.Range("K65").Value = 16080.8825
$ws.Range("L65").Value = 30642.8565
$ws.Range("M65").Value = -12960.8825
$ws.Range("N65").Value = -36882.85649999999
$ws.Range("H86").Value = 1677.3889
$ws.Range("I86").Value = 1539.9333
$ws.Range("J86").Value = 2364.6667
$ws.Range("K86").Value = 1539.9333
$ws.Range("L86").Value = 2364.6667
$ws.Range("M86").Value = -416.9332999999999
$ws.Range("N86").Value = -4610.6667
$ws.Range("H88").Value = 5557610.5
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 6946638
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 6946638
$ws.Range("M88").Value = -1094
$ws.Range("N88").Value = -6947450
$ws.Range("H89").Value = 1677.3889
$ws.Range("I89").Value = 1539.9333
$ws.Range("J89").Value = 2364.6667
$ws.Range("K89").Value = 7699.666499999999
$ws.Range("L89").Value = 11823.3335
$ws.Range("M89").Value = -2083.666499999999
$ws.Range("N89").Value = -23055.3335
$ws.Range("H91").Value = 5557610.5
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 6946638
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 6946638
$ws.Range("M91").Value = -96
$ws.Range("N91").Value = -6949446
$ws.Range("H96").Value = 342
$ws.Range("I96").Value = 319.42856
$ws.Range("J96").Value = 500
$ws.Range("K96").Value = 958.28568
$ws.Range("L96").Value = 1500
$ws.Range("M96").Value = 414.71432
$ws.Range("N96").Value = -4246
$ws.Range("H107").Value = 370810.3
$ws.Range("I107").Value = 529502.9399999999
$ws.Range("J107").Value = 527.44446
$ws.Range("K107").Value = 529502.9399999999
$ws.Range("L107").Value = 527.44446
$ws.Range("M107").Value = -527582.9399999999
$ws.Range("N107").Value = -4367.44446
$ws.Range("H112").Value = 6199443
$ws.Range("I112").Value = 750
$ws.Range("K112").Value = 2250
$ws.Range("M112").Value = -1142
$ws.Range("H137").Value = 43479944
$ws.Range("I137").Value = 62501056
$ws.Range("J137").Value = 3121.8572
$ws.Range("K137").Value = 187503168
$ws.Range("L137").Value = 9365.571599999999
$ws.Range("M137").Value = -187500618
$ws.Range("N137").Value = -14465.5716

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1238
$ws.Range("I45").Value = 1139.4
$ws.Range("J45").Value = 1566.6666
$ws.Range("K45").Value = 1139.4
$ws.Range("L45").Value = 1566.6666
$ws.Range("M45").Value = -762.4000000000001
$ws.Range("N45").Value = -2320.6666
$ws.Range("H102").Value = 2000
$ws.Range("I102").Value = 2000
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2000
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -378
$ws.Range("H125").Value = 34000
$ws.Range("J125").Value = 34000
$ws.Range("L125").Value = 34000
$ws.Range("N125").Value = -43840

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 4142.857
$ws.Range("J33").Value = 5000
$ws.Range("L33").Value = 5000
$ws.Range("N33").Value = -5672
$ws.Range("H86").Value = 7487.3335
$ws.Range("I86").Value = 1651.1111
$ws.Range("J86").Value = 13323.556
$ws.Range("K86").Value = 1651.1111
$ws.Range("L86").Value = 13323.556
$ws.Range("M86").Value = -528.1111000000001
$ws.Range("N86").Value = -15569.556
$ws.Range("H89").Value = 7487.3335
$ws.Range("I89").Value = 1651.1111
$ws.Range("J89").Value = 13323.556
$ws.Range("K89").Value = 8255.5555
$ws.Range("L89").Value = 66617.78
$ws.Range("M89").Value = -2639.5555
$ws.Range("N89").Value = -77849.78
$ws.Range("H94").Value = 601.76
$ws.Range("I94").Value = 520.2273
$ws.Range("J94").Value = 1199.6666
$ws.Range("K94").Value = 520.2273
$ws.Range("L94").Value = 1199.6666
$ws.Range("M94").Value = -69.22730000000001
$ws.Range("N94").Value = -2101.6666
$ws.Range("H99").Value = 1561.25
$ws.Range("I99").Value = 1098
$ws.Range("K99").Value = 1098
$ws.Range("M99").Value = 400

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1530.75
$ws.Range("I31").Value = 1490.2632
$ws.Range("J31").Value = 2300
$ws.Range("K31").Value = 1490.2632
$ws.Range("L31").Value = 2300
$ws.Range("M31").Value = -1195.2632
$ws.Range("N31").Value = -2890
$ws.Range("H34").Value = 1530.75
$ws.Range("I34").Value = 1490.2632
$ws.Range("J34").Value = 2300
$ws.Range("K34").Value = 1490.2632
$ws.Range("L34").Value = 2300
$ws.Range("M34").Value = -1288.2632
$ws.Range("N34").Value = -2704

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 10417364
$ws.Range("I113").Value = 479.2857
$ws.Range("K113").Value = 1437.8571
$ws.Range("M113").Value = 732.1428999999998

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 23480
$ws.Range("J57").Value = 23480
$ws.Range("L57").Value = 23480
$ws.Range("N57").Value = -25120
$ws.Range("H97").Value = 1298
$ws.Range("I97").Value = 1122.5
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 1122.5
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -626.5
$ws.Range("N97").Value = -2992
$ws.Range("H122").Value = 927752.9399999999
$ws.Range("I122").Value = 1853285.1
$ws.Range("K122").Value = 5559855.300000001
$ws.Range("M122").Value = -5557405.300000001
$ws.Range("H132").Value = 2747.762
$ws.Range("I132").Value = 2719.913
$ws.Range("J132").Value = 2781.4736
$ws.Range("K132").Value = 8159.739
$ws.Range("L132").Value = 8344.4208
$ws.Range("M132").Value = -5629.739
$ws.Range("N132").Value = -13404.4208

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7963.8184
$ws.Range("I61").Value = 9800.308000000001
$ws.Range("J61").Value = 5311.1113
$ws.Range("K61").Value = 9800.308000000001
$ws.Range("L61").Value = 5311.1113
$ws.Range("M61").Value = -9598.308000000001
$ws.Range("N61").Value = -5715.1113
$ws.Range("H93").Value = 1006.2857
$ws.Range("I93").Value = 904.7273
$ws.Range("J93").Value = 1378.6666
$ws.Range("K93").Value = 904.7273
$ws.Range("L93").Value = 1378.6666
$ws.Range("M93").Value = 343.2727
$ws.Range("N93").Value = -3874.6666
$ws.Range("H100").Value = 2376.7334
$ws.Range("I100").Value = 1753.4667
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1753.4667
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1212.4667
$ws.Range("N100").Value = -4082
$ws.Range("H113").Value = 7963.8184
$ws.Range("I113").Value = 9800.308000000001
$ws.Range("J113").Value = 5311.1113
$ws.Range("K113").Value = 9800.308000000001
$ws.Range("L113").Value = 5311.1113
$ws.Range("M113").Value = -7630.308000000001
$ws.Range("N113").Value = -9651.1113
